# Update Behemoth_Profits leve-profit market data snapshot (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 3123.5
$ws.Range("I70").Value = 3123.5
$ws.Range("K70").Value = 9370.5
$ws.Range("M70").Value = -9100.5
# Row 73
$ws.Range("H73").Value = 3123.5
$ws.Range("I73").Value = 3123.5
$ws.Range("K73").Value = 9370.5
$ws.Range("M73").Value = -8434.5
# Row 100
$ws.Range("H100").Value = 1032.8572
$ws.Range("I100").Value = 1004
$ws.Range("K100").Value = 1004
$ws.Range("M100").Value = -463
# Row 138
$ws.Range("H138").Value = 3488.3
$ws.Range("J138").Value = 3559.4111
$ws.Range("L138").Value = 10678.2333
$ws.Range("N138").Value = -20958.2333

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 22791.344
$ws.Range("I32").Value = 10424.303
$ws.Range("K32").Value = 10424.303
$ws.Range("M32").Value = -10137.303
# Row 45
$ws.Range("H45").Value = 29414590
$ws.Range("I45").Value = 33335534
$ws.Range("J45").Value = 7507
$ws.Range("K45").Value = 33335534
$ws.Range("L45").Value = 7507
$ws.Range("M45").Value = -33335157
$ws.Range("N45").Value = -8261
# Row 61
$ws.Range("H61").Value = 6252934
$ws.Range("I61").Value = 2910.3713
$ws.Range("J61").Value = 50003100
$ws.Range("K61").Value = 2910.3713
$ws.Range("L61").Value = 50003100
$ws.Range("M61").Value = -2698.3713
$ws.Range("N61").Value = -50003524
# Row 88
$ws.Range("H88").Value = 3748.1
$ws.Range("I88").Value = 4045.25
$ws.Range("J88").Value = 3550
$ws.Range("K88").Value = 4045.25
$ws.Range("L88").Value = 3550
$ws.Range("M88").Value = -3639.25
$ws.Range("N88").Value = -4362
# Row 91
$ws.Range("H91").Value = 3748.1
$ws.Range("I91").Value = 4045.25
$ws.Range("J91").Value = 3550
$ws.Range("K91").Value = 4045.25
$ws.Range("L91").Value = 3550
$ws.Range("M91").Value = -2641.25
$ws.Range("N91").Value = -6358
# Row 136
$ws.Range("H136").Value = 6252934
$ws.Range("I136").Value = 2910.3713
$ws.Range("J136").Value = 50003100
$ws.Range("K136").Value = 8731.1139
$ws.Range("L136").Value = 150009300
$ws.Range("M136").Value = -6181.1139
$ws.Range("N136").Value = -150014400

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 97
$ws.Range("H97").Value = 29323.166
$ws.Range("I97").Value = 19487.25
$ws.Range("K97").Value = 19487.25
$ws.Range("M97").Value = -18496.25
# Row 110
$ws.Range("H110").Value = 49296.332
$ws.Range("J110").Value = 49296.332
$ws.Range("L110").Value = 49296.332
$ws.Range("N110").Value = -57476.332

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 940197.75
$ws.Range("J31").Value = 2331485.8
$ws.Range("L31").Value = 2331485.8
$ws.Range("N31").Value = -2332075.8
# Row 34
$ws.Range("H34").Value = 940197.75
$ws.Range("J34").Value = 2331485.8
$ws.Range("L34").Value = 2331485.8
$ws.Range("N34").Value = -2331889.8
# Row 88
$ws.Range("H88").Value = 24275.572
$ws.Range("J88").Value = 24275.572
$ws.Range("L88").Value = 24275.572
$ws.Range("N88").Value = -25087.572
# Row 91
$ws.Range("H91").Value = 24275.572
$ws.Range("J91").Value = 24275.572
$ws.Range("L91").Value = 24275.572
$ws.Range("N91").Value = -27083.572
# Row 105
$ws.Range("H105").Value = 908.5
$ws.Range("I105").Value = 809.2727
$ws.Range("K105").Value = 809.2727
$ws.Range("M105").Value = 937.7273
# Row 122
$ws.Range("H122").Value = 2398.923
$ws.Range("I122").Value = 2392.5557
$ws.Range("K122").Value = 7177.6671
$ws.Range("M122").Value = -4727.6671
# Row 132
$ws.Range("H132").Value = 7425.25
$ws.Range("I132").Value = 2982.1667
$ws.Range("J132").Value = 20754.5
$ws.Range("K132").Value = 8946.500100000001
$ws.Range("L132").Value = 62263.5
$ws.Range("M132").Value = -6416.500100000001
$ws.Range("N132").Value = -67323.5
# Row 134
$ws.Range("H134").Value = 2539.3157
$ws.Range("J134").Value = 20014
$ws.Range("L134").Value = 60042
$ws.Range("N134").Value = -65112

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 155.49019
$ws.Range("I2").Value = 77.72727
$ws.Range("J2").Value = 214.48276
$ws.Range("K2").Value = 466.36362
$ws.Range("L2").Value = 1286.89656
$ws.Range("M2").Value = -353.36362
$ws.Range("N2").Value = -1512.89656
# Row 60
$ws.Range("H60").Value = 1927.6428
$ws.Range("I60").Value = 798.7
$ws.Range("K60").Value = 2396.1
$ws.Range("M60").Value = -2145.1
# Row 68
$ws.Range("H68").Value = 2098.818
$ws.Range("I68").Value = 2054.5
$ws.Range("K68").Value = 6163.5
$ws.Range("M68").Value = -5352.5
# Row 71
$ws.Range("H71").Value = 2098.818
$ws.Range("I71").Value = 2054.5
$ws.Range("K71").Value = 18490.5
$ws.Range("M71").Value = -14434.5
# Row 76
$ws.Range("H76").Value = 7000
$ws.Range("J76").Value = 7000
$ws.Range("L76").Value = 21000
$ws.Range("N76").Value = -21766
# Row 79
$ws.Range("H79").Value = 7000
$ws.Range("J79").Value = 7000
$ws.Range("L79").Value = 21000
$ws.Range("N79").Value = -23652
# Row 86
$ws.Range("H86").Value = 532.7692
$ws.Range("I86").Value = 490.85715
$ws.Range("J86").Value = 581.6667
$ws.Range("K86").Value = 1472.57145
$ws.Range("L86").Value = 1745.0001
$ws.Range("M86").Value = -286.5714499999999
$ws.Range("N86").Value = -4117.0001
# Row 87
$ws.Range("H87").Value = 22525
$ws.Range("I87").Value = 100
$ws.Range("K87").Value = 300
$ws.Range("M87").Value = 948
# Row 88
$ws.Range("H88").Value = 4750
$ws.Range("I88").Value = 3500
$ws.Range("J88").Value = 4863.636
$ws.Range("K88").Value = 10500
$ws.Range("L88").Value = 14590.908
$ws.Range("M88").Value = -10072
$ws.Range("N88").Value = -15446.908
# Row 89
$ws.Range("H89").Value = 532.7692
$ws.Range("I89").Value = 490.85715
$ws.Range("J89").Value = 581.6667
$ws.Range("K89").Value = 4417.71435
$ws.Range("L89").Value = 5235.0003
$ws.Range("M89").Value = 1510.28565
$ws.Range("N89").Value = -17091.0003
# Row 90
$ws.Range("H90").Value = 22525
$ws.Range("I90").Value = 100
$ws.Range("K90").Value = 900
$ws.Range("M90").Value = 5340
# Row 91
$ws.Range("H91").Value = 4750
$ws.Range("I91").Value = 3500
$ws.Range("J91").Value = 4863.636
$ws.Range("K91").Value = 10500
$ws.Range("L91").Value = 14590.908
$ws.Range("M91").Value = -9018
$ws.Range("N91").Value = -17554.908
# Row 115
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
# Row 131
$ws.Range("H131").Value = 14548.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 14548.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 43645.5
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -53725.5
# Row 132
$ws.Range("H132").Value = 1531.3125
$ws.Range("J132").Value = 1383.25
$ws.Range("L132").Value = 12449.25
$ws.Range("N132").Value = -17509.25

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3162.8572
$ws.Range("I46").Value = 1860
$ws.Range("J46").Value = 4900
$ws.Range("K46").Value = 1860
$ws.Range("L46").Value = 4900
$ws.Range("M46").Value = -1672
$ws.Range("N46").Value = -5276
# Row 93
$ws.Range("H93").Value = 1790.7333
$ws.Range("I93").Value = 1525.65
$ws.Range("J93").Value = 2320.9
$ws.Range("K93").Value = 1525.65
$ws.Range("L93").Value = 2320.9
$ws.Range("M93").Value = -277.6500000000001
$ws.Range("N93").Value = -4816.9
# Row 132
$ws.Range("H132").Value = 6146.067
$ws.Range("I132").Value = 7157.5835
$ws.Range("J132").Value = 2100
$ws.Range("K132").Value = 21472.7505
$ws.Range("L132").Value = 6300
$ws.Range("M132").Value = -18942.7505
$ws.Range("N132").Value = -11360

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 2890.1
$ws.Range("I136").Value = 3067.3
$ws.Range("J136").Value = 2535.7
$ws.Range("K136").Value = 9201.900000000001
$ws.Range("L136").Value = 7607.099999999999
$ws.Range("M136").Value = -6651.900000000001
$ws.Range("N136").Value = -12707.1
